$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Value = "MinWeight"
$ws.Range("C1").Value = "MaxWeight"

# Ticker data rows (switched from arithmetic to log returns -> refreshed tickers/weights)
$ws.Range("A2").Value = "FDX"
$ws.Range("B2").Value = 0.01
$ws.Range("C2").Value = 1
$ws.Range("A3").Value = "WMT"
$ws.Range("B3").Value = 0.01
$ws.Range("C3").Value = 1
$ws.Range("A4").Value = "AMZN"
$ws.Range("B4").Value = 0.01
$ws.Range("C4").Value = 1
$ws.Range("A5").Value = "CORE"
$ws.Range("B5").Value = 0.01
$ws.Range("C5").Value = 1
$ws.Range("A6").Value = "BLK"
$ws.Range("B6").Value = 0.01
$ws.Range("C6").Value = 1
$ws.Range("A7").Value = "LMT"
$ws.Range("B7").Value = 0.01
$ws.Range("C7").Value = 1
$ws.Range("A8").Value = "ORCL"
$ws.Range("B8").Value = 0.01
$ws.Range("C8").Value = 1
$ws.Range("A9").Value = "NTRS"
$ws.Range("B9").Value = 0.01
$ws.Range("C9").Value = 1
$ws.Range("A10").Value = "TSM"
$ws.Range("B10").Value = 0.01
$ws.Range("C10").Value = 1
$ws.Range("A11").Value = "SJM"
$ws.Range("B11").Value = 0.01
$ws.Range("C11").Value = 1
$ws.Range("A12").Value = "MDLZ"
$ws.Range("B12").Value = 0.01
$ws.Range("C12").Value = 1
$ws.Range("A13").Value = "REGI"
$ws.Range("B13").Value = 0.01
$ws.Range("C13").Value = 1
$ws.Range("A14").Value = "V"
$ws.Range("B14").Value = 0.01
$ws.Range("C14").Value = 1
$ws.Range("A15").Value = "MSFT"
$ws.Range("B15").Value = 0.01
$ws.Range("C15").Value = 1
$ws.Range("A16").Value = "JNJ"
$ws.Range("B16").Value = 0.01
$ws.Range("C16").Value = 1
$ws.Range("A17").Value = "TPH"
$ws.Range("B17").Value = 0.01
$ws.Range("C17").Value = 1
$ws.Range("A18").Value = "NRZ"
$ws.Range("B18").Value = 0.01
$ws.Range("C18").Value = 1
$ws.Range("A19").Value = "VIRT"
$ws.Range("B19").Value = 0.01
$ws.Range("C19").Value = 1
$ws.Range("A20").Value = "AXP"
$ws.Range("B20").Value = 0.01
$ws.Range("C20").Value = 1
$ws.Range("A21").Value = "BX"
$ws.Range("B21").Value = 0.01
$ws.Range("C21").Value = 1
$ws.Range("A22").Value = "CNC"
$ws.Range("B22").Value = 0.01
$ws.Range("C22").Value = 1
$ws.Range("A23").Value = "LDOS"
$ws.Range("B23").Value = 0.01
$ws.Range("C23").Value = 1
$ws.Range("A24").Value = "MDT"
$ws.Range("B24").Value = 0.01
$ws.Range("C24").Value = 1
$ws.Range("A25").Value = "MRK"
$ws.Range("B25").Value = 0.01
$ws.Range("C25").Value = 1
$ws.Range("A26").Value = "NKE"
$ws.Range("B26").Value = 0.01
$ws.Range("C26").Value = 1
$ws.Range("A27").Value = "OHI"
$ws.Range("B27").Value = 0.01
$ws.Range("C27").Value = 1
$ws.Range("A28").Value = "PFE"
$ws.Range("B28").Value = 0.01
$ws.Range("C28").Value = 1
$ws.Range("A29").Value = "SYY"
$ws.Range("B29").Value = 0.01
$ws.Range("C29").Value = 1
$ws.Range("A30").Value = "GLTR"
$ws.Range("B30").Value = 0.01
$ws.Range("C30").Value = 1
$ws.Range("A31").Value = "EMB"
$ws.Range("B31").Value = 0.01
$ws.Range("C31").Value = 1
$ws.Range("A32").Value = "TLT"
$ws.Range("B32").Value = 0.01
$ws.Range("C32").Value = 1

# Header comments
$ws.Range("B1").AddComment("Sum of column must be less than 1") | Out-Null
$ws.Range("C1").AddComment("Sum of column must be greater than 1") | Out-Null

# Selection moves to B4
$ws.Range("B4").Select() | Out-Null
